# Minor upgrades in documentation (input, engine and index).
# Delete the obsolete "Weight distribution" row (1961/Global) from the
# "Input data" sheet, which shifts all subsequent rows up by one.

$wb = $excel.ActiveWorkbook

# Work on the "Input data" sheet and make it the active one.
$ws = $wb.Worksheets.Item("Input data")
$ws.Activate()

# Row 136 currently holds: France | 1961 | Global | Weight distribution | 0.05
# Select the full row (as a user would by clicking the row header) and
# delete it entirely, shifting the rows below upward.
$row = $ws.Rows.Item(136)
$row.Select()
$row.Delete()

# Scroll/position the view similarly to how it ends up after the edit:
# frozen header row, top-left visible cell at A130, and the newly
# shifted row 136 selected as a whole row.
$ws.Application.ActiveWindow.ScrollRow = 130
$ws.Range("A136:XFD136").Select()

# Make sure "Input data" remains the active/selected sheet when saved.
$ws.Select()
